$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 161.75
$ws.Range("J19").Value = 150
$ws.Range("L19").Value = 150
$ws.Range("N19").Value = -500
$ws.Range("H33").Value = 485.06668
$ws.Range("I33").Value = 173
$ws.Range("K33").Value = 173
$ws.Range("M33").Value = 56
$ws.Range("H40").Value = 3157.6924
$ws.Range("I40").Value = 1599
$ws.Range("J40").Value = 3287.5833
$ws.Range("K40").Value = 1599
$ws.Range("L40").Value = 3287.5833
$ws.Range("M40").Value = -1424
$ws.Range("N40").Value = -3637.5833
$ws.Range("H58").Value = 2022.6923
$ws.Range("I58").Value = 957.5
$ws.Range("J58").Value = 2216.3635
$ws.Range("K58").Value = 2872.5
$ws.Range("L58").Value = 6649.0905
$ws.Range("M58").Value = -2722.5
$ws.Range("N58").Value = -6949.0905
$ws.Range("H62").Value = 4742.6665
$ws.Range("I62").Value = 4767.143
$ws.Range("J62").Value = 4400
$ws.Range("K62").Value = 4767.143
$ws.Range("L62").Value = 4400
$ws.Range("M62").Value = -4143.143
$ws.Range("N62").Value = -5648
$ws.Range("H65").Value = 4742.6665
$ws.Range("I65").Value = 4767.143
$ws.Range("J65").Value = 4400
$ws.Range("K65").Value = 23835.715
$ws.Range("L65").Value = 22000
$ws.Range("M65").Value = -20715.715
$ws.Range("N65").Value = -28240
$ws.Range("H95").Value = 15062
$ws.Range("J95").Value = 15062
$ws.Range("L95").Value = 15062
$ws.Range("N95").Value = -20554
$ws.Range("H98").Value = 1938.3125
$ws.Range("I98").Value = 643.1111
$ws.Range("K98").Value = 643.1111
$ws.Range("M98").Value = 854.8889
$ws.Range("H122").Value = 1938.3125
$ws.Range("I122").Value = 643.1111
$ws.Range("K122").Value = 1929.3333
$ws.Range("M122").Value = 520.6667000000002
$ws.Range("H132").Value = 29415096
$ws.Range("I132").Value = 34485990
$ws.Range("K132").Value = 103457970
$ws.Range("M132").Value = -103455440
$ws.Range("H135").Value = 811.75
$ws.Range("J135").Value = 2999.5
$ws.Range("L135").Value = 26995.5
$ws.Range("N135").Value = -32065.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H30").Value = 561.8
$ws.Range("I30").Value = 452.25
$ws.Range("J30").Value = 1000
$ws.Range("K30").Value = 452.25
$ws.Range("L30").Value = 1000
$ws.Range("M30").Value = -302.25
$ws.Range("N30").Value = -1300
$ws.Range("H61").Value = 2999
$ws.Range("I61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("M61").ClearContents()
$ws.Range("H97").Value = 1146.4286
$ws.Range("I97").Value = 846.6667
$ws.Range("J97").Value = 2945
$ws.Range("K97").Value = 846.6667
$ws.Range("L97").Value = 2945
$ws.Range("M97").Value = -350.6667
$ws.Range("N97").Value = -3937
$ws.Range("H98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("L98").ClearContents()
$ws.Range("N98").Value = 0
$ws.Range("H105").Value = 51000
$ws.Range("J105").Value = 51000
$ws.Range("L105").Value = 51000
$ws.Range("N105").Value = -57988
$ws.Range("H110").Value = 4301
$ws.Range("I110").Value = 1505
$ws.Range("K110").Value = 1505
$ws.Range("M110").Value = 540
$ws.Range("H132").Value = 4012.6667
$ws.Range("I132").Value = 3708.8696
$ws.Range("K132").Value = 11126.6088
$ws.Range("M132").Value = -8596.6088
$ws.Range("H136").Value = 2999
$ws.Range("I136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("M136").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 213.42857
$ws.Range("I80").Value = 76.333336
$ws.Range("K80").Value = 76.333336
$ws.Range("M80").Value = 921.666664
$ws.Range("H82").Value = 19145.6
$ws.Range("H83").Value = 213.42857
$ws.Range("I83").Value = 76.333336
$ws.Range("K83").Value = 381.66668
$ws.Range("M83").Value = 4610.33332
$ws.Range("H85").Value = 19145.6
$ws.Range("H100").Value = 41938.332
$ws.Range("J100").Value = 41938.332
$ws.Range("L100").Value = 41938.332
$ws.Range("N100").Value = -44102.332
$ws.Range("H107").Value = 992.1539
$ws.Range("I107").Value = 616.55554
$ws.Range("J107").Value = 1837.25
$ws.Range("K107").Value = 616.55554
$ws.Range("L107").Value = 1837.25
$ws.Range("M107").Value = 1303.44446
$ws.Range("N107").Value = -5677.25
$ws.Range("H127").Value = 35998
$ws.Range("J127").Value = 35998
$ws.Range("L127").Value = 35998
$ws.Range("N127").Value = -45918

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 3376.5
$ws.Range("I58").Value = 2314.75
$ws.Range("K58").Value = 2314.75
$ws.Range("M58").Value = -2111.75
$ws.Range("H122").Value = 3812.9285
$ws.Range("I122").Value = 4431.3335
$ws.Range("J122").Value = 2699.8
$ws.Range("K122").Value = 13294.0005
$ws.Range("L122").Value = 8099.400000000001
$ws.Range("M122").Value = -10844.0005
$ws.Range("N122").Value = -12999.4
$ws.Range("H136").Value = 3376.5
$ws.Range("I136").Value = 2314.75
$ws.Range("K136").Value = 6944.25
$ws.Range("M136").Value = -4394.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 2603.4167
$ws.Range("I129").Value = 1750
$ws.Range("J129").Value = 3030.125
$ws.Range("K129").Value = 5250
$ws.Range("L129").Value = 9090.375
$ws.Range("M129").Value = -250
$ws.Range("N129").Value = -19090.375
$ws.Range("H140").Value = 2000
$ws.Range("I140").Value = 2000
$ws.Range("K140").Value = 6000
$ws.Range("M140").Value = -820

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H92").Value = 3049.3333
$ws.Range("J92").Value = 3049.3333
$ws.Range("L92").Value = 3049.3333
$ws.Range("N92").Value = -6793.3333
$ws.Range("H132").Value = 982
$ws.Range("I132").Value = 982
$ws.Range("K132").Value = 2946
$ws.Range("M132").Value = -416

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H97").Value = 6000.5
$ws.Range("J97").Value = 6000.5
$ws.Range("L97").Value = 6000.5
$ws.Range("N97").Value = -7982.5
$ws.Range("H103").Value = 35466.332
$ws.Range("J103").Value = 35466.332
$ws.Range("L103").Value = 35466.332
$ws.Range("N103").Value = -37810.332
$ws.Range("H106").Value = 11107.375
$ws.Range("J106").Value = 11107.375
$ws.Range("L106").Value = 11107.375
$ws.Range("N106").Value = -13631.375
$ws.Range("H114").Value = 41999.5
$ws.Range("J114").Value = 41999.5
$ws.Range("L114").Value = 41999.5
$ws.Range("N114").Value = -50677.5
$ws.Range("H122").Value = 3531.1667
$ws.Range("J122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("N122").ClearContents()
$ws.Range("H136").Value = 3212.2856
$ws.Range("I136").Value = 3212.2856
$ws.Range("K136").Value = 9636.856800000001
$ws.Range("M136").Value = -7086.856800000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 1500
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").ClearContents()
$ws.Range("H65").Value = 1500
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").ClearContents()
$ws.Range("H76").Value = 60000
$ws.Range("J76").Value = 60000
$ws.Range("L76").Value = 60000
$ws.Range("N76").Value = -60630
$ws.Range("H79").Value = 60000
$ws.Range("J79").Value = 60000
$ws.Range("L79").Value = 60000
$ws.Range("N79").Value = -62184
$ws.Range("H80").Value = 34867.332
$ws.Range("J80").Value = 34867.332
$ws.Range("L80").Value = 34867.332
$ws.Range("N80").Value = -36863.332
$ws.Range("H82").Value = 0
$ws.Range("J82").Value = 0
$ws.Range("L82").ClearContents()
$ws.Range("N82").Value = 0
$ws.Range("H83").Value = 34867.332
$ws.Range("J83").Value = 34867.332
$ws.Range("L83").Value = 104601.996
$ws.Range("N83").Value = -114585.996
$ws.Range("H85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("L85").ClearContents()
$ws.Range("N85").Value = 0
$ws.Range("H95").Value = 14672
$ws.Range("J95").Value = 14672
$ws.Range("L95").Value = 14672
$ws.Range("N95").Value = -20164
$ws.Range("H104").Value = 26999.75
$ws.Range("J104").Value = 26999.75
$ws.Range("L104").Value = 26999.75
$ws.Range("N104").Value = -33987.75
$ws.Range("H105").Value = 13582.5
$ws.Range("J105").Value = 13582.5
$ws.Range("L105").Value = 13582.5
$ws.Range("N105").Value = -20570.5
$ws.Range("H111").Value = 39660.75
$ws.Range("J111").Value = 39660.75
$ws.Range("L111").Value = 39660.75
$ws.Range("N111").Value = -47840.75
